# Apply the scheduled text replacements (date header + multiplication
# problems) described by the commit diff.
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2024-12-13 Friday"; New = "2024-12-14 Saturday" },
    @{ Old = "411×2=822";         New = "266×7=1862" },
    @{ Old = "409×2=818";         New = "423×5=2115" },
    @{ Old = "625×2=1250";        New = "927×4=3708" },
    @{ Old = "288×2=576";         New = "729×8=5832" },
    @{ Old = "862×6=5172";        New = "788×9=7092" },
    @{ Old = "231×9=2079";        New = "561×6=3366" },
    @{ Old = "266×6=1596";        New = "820×5=4100" },
    @{ Old = "149×4=596";         New = "105×8=840" },
    @{ Old = "292×6=1752";        New = "105×2=210" },
    @{ Old = "644×4=2576";        New = "386×5=1930" },
    @{ Old = "711×8=5688";        New = "933×4=3732" },
    @{ Old = "462×7=3234";        New = "297×3=891" },
    @{ Old = "225×8=1800";        New = "870×2=1740" },
    @{ Old = "995×2=1990";        New = "947×2=1894" },
    @{ Old = "826×8=6608";        New = "341×7=2387" },
    @{ Old = "237×6=1422";        New = "444×9=3996" },
    @{ Old = "332×7=2324";        New = "629×3=1887" },
    @{ Old = "955×4=3820";        New = "755×8=6040" },
    @{ Old = "338×2=676";         New = "688×3=2064" },
    @{ Old = "123×5=615";         New = "583×2=1166" },
    @{ Old = "276×4=1104";        New = "863×2=1726" },
    @{ Old = "345×3=1035";        New = "713×5=3565" },
    @{ Old = "576×3=1728";        New = "192×7=1344" },
    @{ Old = "736×5=3680";        New = "275×3=825" },
    @{ Old = "160×3=480";         New = "531×9=4779" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
